{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst newText =\n  \"As every school in America transitioned to online learning during the COVID-19 lockdowns, I was the technical/development lead on the team who supported all SRE and product engineering teams, working on core platforms and services.\";\n\n// Find the existing first bullet under \"Principal Cloud and Platform Engineer\"\n// \u2014 the \"Authored or edited over 1,800 Confluence documents.\" item \u2014 so the\n// new bullet can be inserted immediately before it. Inserting \"Before\" a\n// paragraph copies that paragraph's formatting (style \"Compact\" + the\n// numId 1002 / ilvl 0 numbering), which matches the other bullets in this\n// job's list.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text || \"\";\n  if (t.indexOf(\"Authored or edited over 1,800\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Authored or edited over 1,800' paragraph\");\n}\n\ntarget.insertParagraph(newText, \"Before\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$newText = \"As every school in America transitioned to online learning during the COVID-19 lockdowns, I was the technical/development lead on the team who supported all SRE and product engineering teams, working on core platforms and services.\"\n\n# Locate the existing first bullet under \"Principal Cloud and Platform Engineer\"\n# (the \"Authored or edited over 1,800 Confluence documents.\" list item) so the\n# new bullet can be inserted immediately before it, inheriting its list /\n# paragraph-style formatting (Compact style, numId 1002, ilvl 0).\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Authored or edited over 1,800*\") {\n        $target = $p\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'Authored or edited over 1,800' paragraph\"\n}\n\n# Insert a new (empty) paragraph before the target; it inherits the target's\n# paragraph formatting (style + numbering), matching the other bullets.\n$target.Range.InsertParagraphBefore()\n\n# The newly inserted paragraph is now at the same index the target used to occupy.\n$newPara = $d.Paragraphs.Item($targetIndex)\n$newPara.Range.Text = $newText\n"}
